$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'41.535.54"
$ws.Range("E2").Value2 = "'  +0.11%  "
$ws.Range("D3").Value2 = "'2.464.09"
$ws.Range("E3").Value2 = "'  -0.22%  "
$ws.Range("D4").Value2 = "'0.999"
$ws.Range("E4").Value2 = "'  -0.72%  "
$ws.Range("D5").Value2 = "'313.85"
$ws.Range("E5").Value2 = "'  +0.60%  "
$ws.Range("D6").Value2 = "'91.13"
$ws.Range("E6").Value2 = "'  -0.50%  "
$ws.Range("D7").Value2 = "'0.549"
$ws.Range("E7").Value2 = "'  +1.50%  "
$ws.Range("E8").Value2 = "'  -0.64%  "
$ws.Range("E9").Value2 = "'  +4.31%  "
$ws.Range("E10").Value2 = "'  -0.54%  "
$ws.Range("D11").Value2 = "'0.0793"
$ws.Range("E11").Value2 = "'  +2.25%  "
$ws.Range("E12").Value2 = "'  +0.73%  "
$ws.Range("D13").Value2 = "'2.846.83"
$ws.Range("E13").Value2 = "'  -0.14%  "
$ws.Range("D14").Value2 = "'6.87"
$ws.Range("E14").Value2 = "'  +0.77%  "
$ws.Range("D15").Value2 = "'15.80"
$ws.Range("E15").Value2 = "'  +3.95%  "
$ws.Range("D16").Value2 = "'2.478.89"
$ws.Range("E16").Value2 = "'  +2.23%  "
$ws.Range("E17").Value2 = "'  -0.15%  "
$ws.Range("D18").Value2 = "'41.544.43"
$ws.Range("E18").Value2 = "'  +0.66%  "
$ws.Range("D19").Value2 = "'6.49"
$ws.Range("E19").Value2 = "'  +3.80%  "
$ws.Range("D20").Value2 = "'0.0₃0940"
$ws.Range("E20").Value2 = "'  +2.73%  "
$ws.Range("D21").Value2 = "'71.05"
$ws.Range("E21").Value2 = "'  +0.69%  "
$ws.Range("D22").Value2 = "'11.23"
$ws.Range("E22").Value2 = "'  +2.91%  "
$ws.Range("D23").Value2 = "'237.88"
$ws.Range("E23").Value2 = "'  +1.31%  "
$ws.Range("E24").Value2 = "'  +0.40%  "
$ws.Range("D25").Value2 = "'1.90"
$ws.Range("E25").Value2 = "'  +2.07%  "
$ws.Range("E26").Value2 = "'  -0.21%  "
$ws.Range("D27").Value2 = "'24.32"
$ws.Range("E27").Value2 = "'  +1.63%  "
$ws.Range("E28").Value2 = "'  +0.29%  "
$ws.Range("D29").Value2 = "'9.67"
$ws.Range("E29").Value2 = "'  +0.14%  "
$ws.Range("D30").Value2 = "'35.29"
$ws.Range("E30").Value2 = "'  -1.03%  "
$ws.Range("D31").Value2 = "'156.21"
$ws.Range("E31").Value2 = "'  +2.67%  "
$ws.Range("D32").Value2 = "'5.43"
$ws.Range("E32").Value2 = "'  +0.53%  "
$ws.Range("E33").Value2 = "'  +0.62%  "
$ws.Range("D34").Value2 = "'0.0756"
$ws.Range("E34").Value2 = "'  +0.73%  "
$ws.Range("D35").Value2 = "'17.12"
$ws.Range("E35").Value2 = "'  -1.01%  "
$ws.Range("D37").Value2 = "'2.86"
$ws.Range("E37").Value2 = "'  -3.53%  "
$ws.Range("E38").Value2 = "'  +1.59%  "
$ws.Range("E39").Value2 = "'  +2.84%  "
$ws.Range("D40").Value2 = "'1.77"
$ws.Range("E40").Value2 = "'  -3.00%  "
$ws.Range("D41").Value2 = "'3.99"
$ws.Range("E41").Value2 = "'  -0.46%  "
$ws.Range("E42").Value2 = "'  -0.94%  "
$ws.Range("D43").Value2 = "'1.958.10"
$ws.Range("E43").Value2 = "'  -0.15%  "
$ws.Range("E44").Value2 = "'  +0.75%  "
$ws.Range("D45").Value2 = "'18.62"
$ws.Range("E45").Value2 = "'  -3.26%  "
$ws.Range("E46").Value2 = "'  -0.63%  "
$ws.Range("D47").Value2 = "'9.01"
$ws.Range("E47").Value2 = "'  +4.83%  "
$ws.Range("D48").Value2 = "'2.705.76"
$ws.Range("E48").Value2 = "'  -0.52%  "
$ws.Range("D49").Value2 = "'96.86"
$ws.Range("E49").Value2 = "'  +1.45%  "
$ws.Range("D50").Value2 = "'67.07"
$ws.Range("E50").Value2 = "'  -0.80%  "
$ws.Range("D51").Value2 = "'0.171"
$ws.Range("E51").Value2 = "'  -1.50%  "
